$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in "Nada" for the third standup block (row 10-12), and the new
# note about what opposes the goal for Paula Andrea (row 11 -> column D)
$ws.Range("C10").Value = "Nada"
$ws.Range("D10").Value = "Nada"

$ws.Range("C11").Value = "Nada"
$ws.Range("D11").Value = "Leer sobre métodos, prácticas y núcleos. Pensar en clientes potenciales."

$ws.Range("C12").Value = "Nada"
$ws.Range("D12").Value = "Nada"

# Update the selected cell to reflect where the author left off editing
$ws.Range("D11").Select()
